# "changing the database for rectifying errors"
#
# The G-column formula built an INSERT statement for [dbo].[slot details]
# that erroneously included the row-id (column A) as the first value in
# the VALUES(...) list. Rectify it by dropping that argument from the
# CONCATENATE() call for every data row (1-200); the shared-formula
# engine will recreate the per-row relative references (B/C/D/E/F) the
# same way a fill-down in the UI would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFormula = '=CONCATENATE("Insert into [dbo].[slot details] values(",B1,",' + "'" + '",C1,"' + "'" + ',",D1,",",E1,",",F1,")")'

$ws.Range("G1:G200").Formula = $newFormula

# Recall the scroll position / selection the author left the sheet in.
$win = $excel.ActiveWindow
$win.ScrollRow = 183
$win.ScrollColumn = 1

$ws.Range("G1:G200").Select()
